$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "36÷4=9, 0"
$t.Cell(1,2).Range.Text = "21÷7=3, 0"
$t.Cell(1,3).Range.Text = "16÷5=3, 1"
$t.Cell(1,4).Range.Text = "13÷6=2, 1"
$t.Cell(1,5).Range.Text = "23÷2=11, 1"

$t.Cell(5,1).Range.Text = "61÷9=6, 7"
$t.Cell(5,2).Range.Text = "30÷4=7, 2"
$t.Cell(5,3).Range.Text = "62÷4=15, 2"
$t.Cell(5,4).Range.Text = "12÷3=4, 0"
$t.Cell(5,5).Range.Text = "14÷3=4, 2"

$t.Cell(9,1).Range.Text = "47÷3=15, 2"
$t.Cell(9,2).Range.Text = "58÷2=29, 0"
$t.Cell(9,3).Range.Text = "14÷6=2, 2"
$t.Cell(9,4).Range.Text = "27÷2=13, 1"
$t.Cell(9,5).Range.Text = "86÷6=14, 2"

$t.Cell(13,1).Range.Text = "91÷8=11, 3"
$t.Cell(13,2).Range.Text = "14÷9=1, 5"
$t.Cell(13,3).Range.Text = "73÷3=24, 1"
$t.Cell(13,4).Range.Text = "67÷7=9, 4"
$t.Cell(13,5).Range.Text = "94÷7=13, 3"

$t.Cell(17,1).Range.Text = "89÷4=22, 1"
$t.Cell(17,2).Range.Text = "61÷9=6, 7"
$t.Cell(17,3).Range.Text = "53÷2=26, 1"
$t.Cell(17,4).Range.Text = "51÷6=8, 3"
$t.Cell(17,5).Range.Text = "60÷4=15, 0"
